$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 8 (shifts old rows 8-11 down to 9-12, duplicating formatting of old row 8)
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new scenario's data
$ws.Range("A8").Value = "CW3M"
$ws.Range("B8").Value = "Baseline C68+ 2010 10/17/20"
$ws.Range("C8").Value = 2010
$ws.Range("D8").Value = 1044.2558590000001
$ws.Range("E8").Value = 1990.4676509999999
$ws.Range("F8").Value = 1.255063
$ws.Range("G8").Value = 327.58108499999997
$ws.Range("H8").Value = 10.610913999999999
$ws.Range("I8").Value = 8.8404570000000007
$ws.Range("J8").Value = 814.39868200000001
$ws.Range("K8").Value = 93.229797000000005
$ws.Range("L8").Value = 1291.7857670000001
$ws.Range("M8").Value = 1165.4420170000001
$ws.Range("N8").Value = 7166.0351559999999
$ws.Range("O8").Value = 29450.638672000001
$ws.Range("P8").Value = -0.473854
$ws.Range("Q8").Value = -0.00014
$ws.Range("R8").Value = 2010

# P8 keeps the default (General) number format, not the inherited one
$ws.Range("P8").ClearFormats()
$ws.Range("P8").Value = -0.473854

# The old (now-shifted) row 8 becomes row 9, which is blank except for Q9 -
# clear out its leftover inherited formatting entirely
$ws.Range("D9:P9").Clear()

# Give P7 a new, more precise number format (0.00000)
$ws.Range("P7").NumberFormat = "0.00000"

# Update selection to match the post-edit state
$ws.Range("N8:O8").Select() | Out-Null
